$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows that only get their B column bumped by +4 ---
$simpleRows = @(38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,59,62,63)
foreach ($r in $simpleRows) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = $cell.Value2 + 4
}

# --- Rows 57 and 58 swap their whole row content (incl. B), then B gets +4 ---
$cols = @(1,2,4,5,6,7,8,17,18,26,28,29,36,37,41)  # A,B,D,E,F,G,H,Q,R,Z,AB,AC,AJ,AK,AO

$row57 = @{}
$row58 = @{}
foreach ($c in $cols) {
    $row57[$c] = $ws.Cells.Item(57, $c).Value2
    $row58[$c] = $ws.Cells.Item(58, $c).Value2
}

foreach ($c in $cols) {
    $ws.Cells.Item(57, $c).Value2 = $row58[$c]
    $ws.Cells.Item(58, $c).Value2 = $row57[$c]
}
$ws.Cells.Item(57, 2).Value2 = $ws.Cells.Item(57, 2).Value2 + 4
$ws.Cells.Item(58, 2).Value2 = $ws.Cells.Item(58, 2).Value2 + 4

# --- Rows 60 and 61 swap their whole row content (incl. B), then B gets +4 ---
$row60 = @{}
$row61 = @{}
foreach ($c in $cols) {
    $row60[$c] = $ws.Cells.Item(60, $c).Value2
    $row61[$c] = $ws.Cells.Item(61, $c).Value2
}

foreach ($c in $cols) {
    $ws.Cells.Item(60, $c).Value2 = $row61[$c]
    $ws.Cells.Item(61, $c).Value2 = $row60[$c]
}
$ws.Cells.Item(60, 2).Value2 = $ws.Cells.Item(60, 2).Value2 + 4
$ws.Cells.Item(61, 2).Value2 = $ws.Cells.Item(61, 2).Value2 + 4
